# Add a new "Umbilical Cable" sub-system row to the On-Site and Inspections
# sheets, duplicating the existing "Export Cable" row's values/format and
# only changing the label - reflecting fewer distinct variables being
# tracked (moorings & installation examples are relied on instead).

$wb = $excel.ActiveWorkbook

# ---- On-Site sheet: new row 5 "Umbilical Cable" (copy of row 4) ----
$ws1 = $wb.Worksheets.Item("On-Site")

# Copy row 4's formatting (including the trailing, value-less R:T cells)
# into row 5 first, ...
$ws1.Range("A4:T4").Copy()
$ws1.Range("A5:T5").PasteSpecial(-4122)

# ... then copy row 4's values into row 5.
$ws1.Range("A4:Q4").Copy()
$ws1.Range("A5:Q5").PasteSpecial(-4163)

# Re-label the duplicated row as "Umbilical Cable".
$ws1.Range("A5").Value = "Umbilical Cable"

$ws1.Activate()
$ws1.Range("A5").Select()

# ---- Inspections sheet: new row 5 "Umbilical Cable" (copy of row 4) ----
$ws3 = $wb.Worksheets.Item("Inspections")

$ws3.Range("A4:O4").Copy()
$ws3.Range("A5:O5").PasteSpecial(-4122)

$ws3.Range("A4:O4").Copy()
$ws3.Range("A5:O5").PasteSpecial(-4163)

$ws3.Range("A5").Value = "Umbilical Cable"

$ws3.Activate()
$ws3.Range("A23").Select()

# Leave the On-Site sheet as the active / selected sheet, matching the
# original workbook state (tabSelected="1" on the On-Site sheet).
$ws1.Activate()
